$wb = $excel.ActiveWorkbook

# --- Sheet "Tableless": fix DateTime values/format (bug fix) ---
$ws = $wb.Worksheets.Item("Tableless")

# Make this sheet the active one (also updates workbook.xml activeTab and
# moves tabSelected from the previously active sheet to this one).
$ws.Activate()

# C3: 45214 (date only) -> 45214.75 (date + time), with a date+time number format
$c3 = $ws.Range("C3")
$c3.NumberFormat = "d/m/yy\ h:mm;@"
$c3.Value = 45214.75

# C4: 45215 (date only) -> 45215.8333333333 (date + time), same format
$c4 = $ws.Range("C4")
$c4.NumberFormat = "d/m/yy\ h:mm;@"
$c4.Value = 45215.8333333333

# Column C needs to widen (to fit the longer date+time text) and be split
# off from column D (which keeps its original width).
$ws.Columns.Item(3).ColumnWidth = 12.5

# Update the visible selection to match the edited cells.
$ws.Range("C3:C4").Select()
